$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 593.075
$ws.Range("I80").Value = 638
$ws.Range("J80").Value = 488.25
$ws.Range("K80").Value = 1914
$ws.Range("L80").Value = 1464.75
$ws.Range("M80").Value = -916
$ws.Range("N80").Value = -3460.75

$ws.Range("H83").Value = 593.075
$ws.Range("I83").Value = 638
$ws.Range("J83").Value = 488.25
$ws.Range("K83").Value = 5742
$ws.Range("L83").Value = 4394.25
$ws.Range("M83").Value = -750
$ws.Range("N83").Value = -14378.25

$ws.Range("H132").Value = 2021.5
$ws.Range("I132").Value = 2061.6924
$ws.Range("K132").Value = 6185.0772
$ws.Range("M132").Value = -3655.0772

$ws.Range("H137").Value = 14376.368
$ws.Range("I137").Value = 4575.6665
$ws.Range("K137").Value = 13726.9995
$ws.Range("M137").Value = -11176.9995

$ws.Range("H138").Value = 4713.479
$ws.Range("J138").Value = 4787.978
$ws.Range("L138").Value = 14363.934
$ws.Range("N138").Value = -24643.934

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 60001
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H20").Value = 60001
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H32").Value = 61844.203
$ws.Range("I32").Value = 19608.268
$ws.Range("K32").Value = 19608.268
$ws.Range("M32").Value = -19321.268

$ws.Range("H61").Value = 4625.9165
$ws.Range("I61").Value = 2579.4736
$ws.Range("K61").Value = 2579.4736
$ws.Range("M61").Value = -2367.4736

$ws.Range("H74").Value = 2593.318
$ws.Range("I74").Value = 1397
$ws.Range("K74").Value = 1397
$ws.Range("M74").Value = -523

$ws.Range("H77").Value = 2593.318
$ws.Range("I77").Value = 1397
$ws.Range("K77").Value = 6985
$ws.Range("M77").Value = -2617

$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

$ws.Range("H136").Value = 4625.9165
$ws.Range("I136").Value = 2579.4736
$ws.Range("K136").Value = 7738.4208
$ws.Range("M136").Value = -5188.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 93558.13
$ws.Range("J99").Value = 1000499.5
$ws.Range("L99").Value = 1000499.5
$ws.Range("N99").Value = -1003495.5

$ws.Range("H134").Value = 11910201
$ws.Range("I134").Value = 5213885.5
$ws.Range("K134").Value = 15641656.5
$ws.Range("M134").Value = -15639121.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12991455
$ws.Range("I31").Value = 22729796
$ws.Range("K31").Value = 22729796
$ws.Range("M31").Value = -22729501

$ws.Range("H34").Value = 12991455
$ws.Range("I34").Value = 22729796
$ws.Range("K34").Value = 22729796
$ws.Range("M34").Value = -22729594

$ws.Range("H132").Value = 4312.5884
$ws.Range("I132").Value = 3736.7856
$ws.Range("K132").Value = 11210.3568
$ws.Range("M132").Value = -8680.356800000001

$ws.Range("H141").Value = 229210.19
$ws.Range("J141").Value = 237987.28
$ws.Range("L141").Value = 237987.28
$ws.Range("N141").Value = -248347.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6743.1665
$ws.Range("I56").Value = 6743.1665
$ws.Range("K56").Value = 6743.1665
$ws.Range("M56").Value = -6213.1665

$ws.Range("H113").Value = 1873.4
$ws.Range("J113").Value = 1873.4
$ws.Range("L113").Value = 5620.200000000001
$ws.Range("N113").Value = -9960.200000000001

$ws.Range("H132").Value = 5730.0557
$ws.Range("I132").Value = 1878.25
$ws.Range("J132").Value = 13433.667
$ws.Range("K132").Value = 16904.25
$ws.Range("L132").Value = 120903.003
$ws.Range("M132").Value = -14374.25
$ws.Range("N132").Value = -125963.003

$ws.Range("H133").Value = 4755.143
$ws.Range("I133").Value = 3381.1667
$ws.Range("K133").Value = 10143.5001
$ws.Range("M133").Value = -5083.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 13259.777
$ws.Range("J26").Value = 13259.777
$ws.Range("L26").Value = 13259.777
$ws.Range("N26").Value = -13819.777

$ws.Range("H47").Value = 30999.5
$ws.Range("J47").Value = 30999.5
$ws.Range("L47").Value = 30999.5
$ws.Range("N47").Value = -32135.5

$ws.Range("H50").Value = 13259.777
$ws.Range("J50").Value = 13259.777
$ws.Range("L50").Value = 13259.777
$ws.Range("N50").Value = -14255.777

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 23274.818
$ws.Range("I22").Value = 675.2
$ws.Range("J22").Value = 42107.832
$ws.Range("K22").Value = 675.2
$ws.Range("L22").Value = 42107.832
$ws.Range("M22").Value = -380.2
$ws.Range("N22").Value = -42697.832

$ws.Range("H27").Value = 23274.818
$ws.Range("I27").Value = 675.2
$ws.Range("J27").Value = 42107.832
$ws.Range("K27").Value = 675.2
$ws.Range("L27").Value = 42107.832
$ws.Range("M27").Value = -568.2
$ws.Range("N27").Value = -42321.832

$ws.Range("H40").Value = 6833.9287
$ws.Range("J40").Value = 7035.5
$ws.Range("L40").Value = 7035.5
$ws.Range("N40").Value = -7307.5

$ws.Range("H46").Value = 1134
$ws.Range("J46").Value = 269
$ws.Range("L46").Value = 269
$ws.Range("N46").Value = -645

$ws.Range("H55").Value = 251.8
$ws.Range("I55").Value = 245.14285
$ws.Range("K55").Value = 245.14285
$ws.Range("M55").Value = -72.14285000000001

$ws.Range("H132").Value = 4300
$ws.Range("I132").Value = 4250
$ws.Range("J132").Value = 4400
$ws.Range("K132").Value = 12750
$ws.Range("L132").Value = 13200
$ws.Range("M132").Value = -10220
$ws.Range("N132").Value = -18260

$ws.Range("H136").Value = 38467800
$ws.Range("I136").Value = 6319.174
$ws.Range("J136").Value = 333339170
$ws.Range("K136").Value = 18957.522
$ws.Range("L136").Value = 1000017510
$ws.Range("M136").Value = -16407.522
$ws.Range("N136").Value = -1000022610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1554.5555
$ws.Range("I136").Value = 1498.875
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4496.625
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1946.625
$ws.Range("N136").Value = -11100
